# Update cryptocurrency price/volume data per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.838.04'
$ws.Range('E2').Value = '  +0.85%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.922.76'
$ws.Range('E3').Value = '  +1.78%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9908'
$ws.Range('E4').Value = '  -1.06%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '250.84'
$ws.Range('E5').Value = '  +2.36%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6036'
$ws.Range('E6').Value = '  +28.05%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9928'
$ws.Range('E7').Value = '  -0.84%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3040'
$ws.Range('E8').Value = '  +4.35%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '24.04'
$ws.Range('E9').Value = '  +7.26%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06660'
$ws.Range('E10').Value = '  +2.56%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7935'
$ws.Range('E11').Value = '  +7.54%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '101.72'
$ws.Range('E12').Value = '  +5.86%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07903'
$ws.Range('E13').Value = '  +1.76%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.892.08'
$ws.Range('E14').Value = '  +0.19%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.331'
$ws.Range('E15').Value = '  +2.74%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '285.59'
$ws.Range('E16').Value = '  +0.58%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.667.23'
$ws.Range('E17').Value = '  +0.03%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.50'
$ws.Range('E18').Value = '  +3.24%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007630'
$ws.Range('E19').Value = '  +1.77%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9966'
$ws.Range('E20').Value = '  -0.38%  '

# Row 21
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.459'
$ws.Range('E21').Value = '  +3.46%  '

# Row 22
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.133.97'
$ws.Range('E22').Value = '  +0.00%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9891'
$ws.Range('E23').Value = '  -1.31%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.583'
$ws.Range('E24').Value = '  +5.19%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.307'
$ws.Range('E25').Value = '  +1.55%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.19'
$ws.Range('E26').Value = '  +0.01%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.51'
$ws.Range('E27').Value = '  +3.33%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.957'
$ws.Range('E28').Value = '  +2.74%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1084'
$ws.Range('E29').Value = '  +11.13%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.353'
$ws.Range('E30').Value = '  +0.26%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.524'
$ws.Range('E31').Value = '  +3.24%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.390'
$ws.Range('E32').Value = '  +2.23%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.272'
$ws.Range('E33').Value = '  +3.36%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04940'
$ws.Range('E34').Value = '  +1.23%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.160'
$ws.Range('E35').Value = '  +2.78%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7153'
$ws.Range('E36').Value = '  +3.25%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.781'
$ws.Range('E37').Value = '  +2.61%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01928'
$ws.Range('E38').Value = '  +1.47%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.905'
$ws.Range('E39').Value = '  +2.39%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '77.14'
$ws.Range('E40').Value = '  +2.17%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.387'
$ws.Range('E41').Value = '  +1.86%  '

# Row 42
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4424'
$ws.Range('E42').Value = '  +3.84%  '

# Row 43
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.010'
$ws.Range('E43').Value = '  +0.32%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8465'
$ws.Range('E44').Value = '  +2.49%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9961'
$ws.Range('E45').Value = '  -0.46%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.09'
$ws.Range('E46').Value = '  +5.90%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.56'
$ws.Range('E47').Value = '  +0.18%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.205'
$ws.Range('E48').Value = '  +3.28%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '35.57'
$ws.Range('E49').Value = '  +0.80%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4076'
$ws.Range('E50').Value = '  +3.41%  '

# Row 51
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '907.24'
$ws.Range('E51').Value = '  -0.22%  '
